$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 85, shifting existing rows 85-213 down to 86-214
$ws.Rows("85:85").Insert()

# Populate the newly inserted row 85 with the new weekly data point
$ws.Cells.Item(85, 1).Value = 8
$ws.Cells.Item(85, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44580
$ws.Cells.Item(85, 5).Value = 4
$ws.Cells.Item(85, 6).Value = 100112012
$ws.Cells.Item(85, 7).Value = "Espinaca"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 3100
$ws.Cells.Item(85, 11).Value = 400
$ws.Cells.Item(85, 12).Value = 500
$ws.Cells.Item(85, 13).Value = 450
$ws.Cells.Item(85, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(85, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(85, 16).Value = 900
$ws.Cells.Item(85, 17).Value = 0.5
$ws.Cells.Item(85, 18).Value = "Hortaliza"
